# SubRES_New_RE_and_Conventional_Trans.xlsx - "Misc" sheet update
# Splits the combined "wind" process/commodity rows into separate
# offshore/onshore rows, renames the ELC_won* wildcard to ELC_wo*,
# extends the fuel list with "bioenergy" and adds a new -life TFM_INS
# override row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")
$ws.Activate()

# --- ~TFM_TOPINS block (rows 40-43 originally) ---------------------------
# Row 41: ElcAgg_Wind now points at the combined wind wildcard ELC_wo*
$ws.Range("D41").Value = "ELC_wo*"

# Row 43: was the combined wind process E[_]W*/wind -> becomes the
# offshore-only row E[_]WOF*/windoff
$ws.Range("C43").Value = "E[_]WOF*"
$ws.Range("D43").Value = "windoff"

# Insert a new row 44 for the onshore counterpart E[_]WON*/windon
$ws.Rows("44:44").Insert()
$ws.Range("C44").Value = "E[_]WON*"
$ws.Range("D44").Value = "windon"
$ws.Range("E44").Value = "IN"

# --- ~TFM_INS block (was rows 46-48, now shifted to 47-49) ---------------
# F49 (was F48): extend the fuel list with bioenergy
$ws.Range("F49").Value = "solar,wind,coal,gas,nuclear,hydro,bioenergy"

# --- New trailing row 50: -life override for coal/gas/nuclear/bioenergy --
$ws.Range("C50").Value = "life"
$ws.Range("D50").Value = 40
$ws.Range("E50").Value = "-life"
$ws.Range("F50").Value = "coal,gas,nuclear,bioenergy"

# E50 needs the same quote-prefix text style as E49 ("-pasti") so the
# leading "-" is kept as literal text formatting, matching the source.
$ws.Range("E49").Copy()
$ws.Range("E50").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cosmetic: restore the selection shown in the saved file -------------
$ws.Range("D42").Select()
